$d = $word.ActiveDocument
$para = $d.Paragraphs.Item(1)
$r = $para.Range

# Set whole text to "G"
$r.Text = "G"

# Get fresh range after paragraph for end, collapse to end of "G"
$para2 = $d.Paragraphs.Item(1)
$rEnd = $para2.Range
$rEnd.Collapse(0)  # wdCollapseEnd = 0
$rEnd.InsertAfter("it")

$para3 = $d.Paragraphs.Item(1)
$rEnd2 = $para3.Range
$rEnd2.Collapse(0)
$rEnd2.InsertAfter(" fichier1")

Write-Output "final: [$($d.Paragraphs.Item(1).Range.Text)]"
